$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.135749
$ws.Range("H2").Value = 0.407247
$ws.Range("M2").Value = 0.789222
$ws.Range("N2").Value = 2.367666
$ws.Range("O2").Value = 0.01341929863527565
$ws.Range("P2").Value = 0.01341929863527565
$ws.Range("Q2").Value = 0.107136097278
$ws.Range("R2").Value = 0.964224875502
$ws.Range("S2").Value = 0.01341929863527565
$ws.Range("T2").Value = 0.01341929863527565

# Row 3
$ws.Range("G3").Value = 0.135749
$ws.Range("H3").Value = 0.407247
$ws.Range("O3").Value = 0.005047365584441773
$ws.Range("P3").Value = 0.005047365584441773
$ws.Range("Q3").Value = 0.040296819152
$ws.Range("R3").Value = 0.362671372368
$ws.Range("S3").Value = 0.005047365584441773
$ws.Range("T3").Value = 0.005047365584441773

# Row 4
$ws.Range("G4").Value = 0.135749
$ws.Range("H4").Value = 0.407247
$ws.Range("M4").Value = 57.61405833333333
$ws.Range("N4").Value = 172.842175
$ws.Range("O4").Value = 0.9796232927683105
$ws.Range("P4").Value = 0.9796232927683105
$ws.Range("Q4").Value = 7.821050804691667
$ws.Range("R4").Value = 70.38945724222501
$ws.Range("S4").Value = 0.9796232927683105
$ws.Range("T4").Value = 0.9796232927683105

# Row 5
$ws.Range("G5").Value = 0.135749
$ws.Range("H5").Value = 0.407247
$ws.Range("M5").Value = 0.1123343333333333
$ws.Range("N5").Value = 0.337003
$ws.Range("O5").Value = 0.001910043011972043
$ws.Range("P5").Value = 0.001910043011972043
$ws.Range("Q5").Value = 0.01524927341566667
$ws.Range("R5").Value = 0.137243460741
$ws.Range("S5").Value = 0.001910043011972043
$ws.Range("T5").Value = 0.001910043011972043
